# Update countries & provincias Spain
# Applies the 10-Jul-2020 07:06 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder two pairs of country names (shared-string order changed) ---
# Italia / Pakistan swap places (rows 14 & 15)
$ws.Range("A14").Value = "Pakistan"
$ws.Range("A15").Value = "Italia"

# Kirguistan moves ahead of Kenia (rows 74, 75, 76 cycle)
$ws.Range("A74").Value = "Kirguistan"
$ws.Range("A75").Value = "Kenia"
$ws.Range("A76").Value = "Noruega"

# --- Refresh case numbers ---

# Brasil (row 5): Casos activos / Recuperados updated
$ws.Range("D5").Value = 1185596
$ws.Range("E5").Value = 504253

# Row 14 now "Pakistan" - new data
$ws.Range("B14").Value = 243599
$ws.Range("C14").Value = 2751
$ws.Range("D14").Value = 149092
$ws.Range("E14").Value = 89449
$ws.Range("G14").Value = 75
$ws.Range("H14").Value = 5058

# Row 15 now "Italia" - new data
$ws.Range("B15").Value = 242363
$ws.Range("D15").Value = 193978
$ws.Range("E15").Value = 13459
$ws.Range("H15").Value = 34926

# Row 73 "Australia" - new data
$ws.Range("B73").Value = 9377
$ws.Range("C73").Value = 318
$ws.Range("E73").Value = 1695

# Row 74 now "Kirguistan" - new data
$ws.Range("B74").Value = 9358
$ws.Range("C74").Value = 511
$ws.Range("D74").Value = 3134
$ws.Range("E74").Value = 6102
$ws.Range("G74").Value = 6
$ws.Range("H74").Value = 122

# Row 75 now "Kenia" - new data
$ws.Range("B75").Value = 8975
$ws.Range("D75").Value = 2657
$ws.Range("E75").Value = 6145
$ws.Range("H75").Value = 173

# Row 76 now "Noruega" - new data
$ws.Range("B76").Value = 8965
$ws.Range("D76").Value = 8138
$ws.Range("E76").Value = 575
$ws.Range("H76").Value = 252

# Row 79 "Venezuela" (country unaffected by the reorder) - new data
$ws.Range("B79").Value = 8372
$ws.Range("D79").Value = 2544
$ws.Range("E79").Value = 5748
$ws.Range("H79").Value = 75

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 07:06"
